$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - count
$ws.Range("B2").Value = 4615394
$ws.Range("C2").Value = 4615394
$ws.Range("D2").Value = 4615394

# Row 3 - mean
$ws.Range("B3").Value = 2631.216873112026
$ws.Range("C3").Value = 50.98690042280242
$ws.Range("D3").Value = 100.5246395323996

# Row 4 - std
$ws.Range("B4").Value = 1687.396185362857
$ws.Range("C4").Value = 0.4607548801340078
$ws.Range("D4").Value = 0.8281038551113546

# Row 5 - min
$ws.Range("B5").Value = -50
$ws.Range("D5").Value = 75

# Row 6 - 25%
$ws.Range("B6").Value = 1520.96

# Row 7 - 50%
$ws.Range("B7").Value = 2083.16
$ws.Range("C7").Value = 50.58

# Row 8 - 75%
$ws.Range("B8").Value = 3639.5375
$ws.Range("D8").Value = 100.72

# Row 9 - max
$ws.Range("B9").Value = 119615.22
$ws.Range("C9").Value = 51.85
$ws.Range("D9").Value = 120.55
